# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.429.42'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.578.93'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.99'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.75'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.575.92'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.17'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.390'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.190.09'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.63'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.580.76'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.969.72'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.03'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.33%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.89'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '393.91'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.584'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.726.18'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.30'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.05'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.79%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +26.55%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.39%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.28%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.578.70'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.03%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.148'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.27'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.59'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.05'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '170.77'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0826'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.844'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.27'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.26'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.50%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.49'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.68'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.49%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.457.27'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0272'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.93%  '
